# Apply the "Dataset complete for TechFest" edit to registration.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Semester values ---
# Rows 9-19: Spring -> Summer
$ws.Range("C9:C19").Value = "Summer"

# Rows 23-32: Summer -> Autumn
$ws.Range("C23:C32").Value = "Autumn"

# --- Remove the trailing rows (46-56), shrinking the dataset down to row 45 ---
$ws.Range("A46:C56").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# --- Update the active selection to match the author's last edit location ---
$ws.Range("D27").Select()

# --- Leftover conditional-format style: applying & then removing a
#     "Duplicate Values" highlight rule registers its dxf (red text /
#     pink fill) in styles.xml without leaving a conditionalFormatting
#     rule bound to the sheet -- matches the captured state exactly.
$fc = $ws.Range("C2:C45").FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 0x06009C
$fc.Interior.Color = 0xCEC7FF
$fc.Delete()

$wb.Save()
